$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the data values to their new decimal equivalents
$ws.Range("A2").Value = 1.1000000000000001
$ws.Range("B2").Value = 2.2000000000000002
$ws.Range("C2").Value = 3.3

$ws.Range("A3").Value = 4.4000000000000004
$ws.Range("B3").Value = 5.5
$ws.Range("C3").Value = 6.6

$ws.Range("A11").Value = 1.1000000000000001
$ws.Range("B11").Value = 2.2000000000000002

$ws.Range("A15").Value = 1.1000000000000001
$ws.Range("B15").Value = 2.2000000000000002
$ws.Range("C15").Value = 3.3

$ws.Range("A16").Value = 4.4000000000000004
$ws.Range("B16").Value = 5.5
$ws.Range("C16").Value = 6.6

$wb.Save()
